$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 / 8: city & address filler text, now wraps ---
$ws.Range("C7").Value = "Recusandae Anim dol "
$ws.Range("C7").WrapText = $true

$ws.Range("C8").Value = "Est quis sunt animi "
$ws.Range("C8").WrapText = $true

# --- Row 11: acknowledgement paragraph, name placeholder changed ---
$ws.Range("B11").Value = "Yo, Ratione voluptatem hRatione voluptatem hRatione voluptatem h certifico haber recibido el carnet de identificación como Servidor de la Empresa y a la vez me responsabilizo por el buen uso del mismo y en caso de la finalización de la relación laboral a la devolución del mismo. `nLibre y voluntariamente me comprometo a depositar en la cuenta de la Empresa Eléctrica Regional Centro Sur C.A.  el valor de 10usd, correspodiente a la reposición por pérdida del carnet de identificación que he recibido. "

# --- Row 17: signature table entry, date now wraps too ---
$ws.Range("B17").Value = "2 de mayo de 1987 "
$ws.Range("B17").WrapText = $true

$ws.Range("C17").Value = "A aut quod voluptas  "

$ws.Range("D17").Value = "Dignissimos voluptat "
$ws.Range("D17").WrapText = $true

$ws.Range("E17").Value = "In tempor exercitati "

# --- Row 22 / 23: responsible person & date ---
$ws.Range("D22").Value = "Non cillum molestiae "

$ws.Range("D23").Value = "2 de mayo de 1987 "
